$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: new "Introduction to DOM" entry ---

# Date (A6) - copy format from A2 (date, bordered) then set value/alignment/format
$ws.Range("A2").Copy($ws.Range("A6"))
$ws.Range("A6").Value = 45286
$ws.Range("A6").HorizontalAlignment = -4108
$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("A6").NumberFormat = "mm-dd-yy"

# Title (B6) - copy format from B2 (bold, centered, bordered)
$ws.Range("B2").Copy($ws.Range("B6"))
$ws.Range("B6").Value = "Introduction to DOM"

# Content (C6) - copy format from C2 (bordered, wrap) then strip horiz/vert alignment to match target (wrap only)
$ws.Range("C2").Copy($ws.Range("C6"))
$ws.Range("C6").Value = "Accessing Dom Elements" + [char]10 + "Looping over list of elements" + [char]10 + "Grab Children/Parent Node(s)" + [char]10 + "Create new DOM Elements" + [char]10 + "Get/Set text to elements" + [char]10 + "Rendering HTML inside elements" + [char]10 + "Add Elements to the DOM" + [char]10 + "Add/Remove/Toggle/Check Classes" + [char]10 + "Event Capturing vs Bubbling"
$ws.Range("C6").HorizontalAlignment = 1
$ws.Range("C6").VerticalAlignment = -4107

# Link (D6) - copy format from D2 (hyperlink style, bordered) then set value + hyperlink
$ws.Range("D2").Copy($ws.Range("D6"))
$ws.Range("D6").Value = "https://scrawny-gouda-144.notion.site/Introduction-to-DOM-12140291d9e94461886732b975092254 "
$ws.Hyperlinks.Add($ws.Range("D6"), "https://scrawny-gouda-144.notion.site/Introduction-to-DOM-12140291d9e94461886732b975092254")

# Row height to match new row (129.6 like the other wrapped content rows)
$ws.Rows.Item(6).RowHeight = 129.6

# --- View state: move selection to F4 ---
$ws.Range("F4").Select()
